# Regenerate orders with updated distance/size codes.
# Mapping (per commit message "regenerate orders with updates distance/sizes"):
#   Distance codes: D80 -> D86, D64 -> D69, D51 -> D55
#   Size codes:      S30 -> S31
# These tokens appear embedded inside many shared strings across the
# Condition / Filename_Left / Filename_Right / Distance / Size columns
# (e.g. "Face12_D80_S25" -> "Face12_D86_S25", "Fixation_D64_l.png" ->
# "Fixation_D69_l.png", the bare "D80"/"D64"/"D51"/"S30" lookup-table
# values, etc.) None of the replacement tokens collide with existing
# text, so a straightforward global substring replace across the used
# range reproduces the diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.UsedRange

$null = $rng.Replace("D80", "D86")
$null = $rng.Replace("D64", "D69")
$null = $rng.Replace("D51", "D55")
$null = $rng.Replace("S30", "S31")
